$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("10:23").Delete()

function Set-Cell($addr, $val, $styleSrc) {
    $ws.Range($addr).Value = $val
    if ($styleSrc) {
        $ws.Range($styleSrc).Copy()
        $ws.Range($addr).PasteSpecial(-4122)
    }
}

Set-Cell "A10" "Objetivos:" "A9"
Set-Cell "B10" "TextB10" "B9"
Set-Cell "C10" "TextC10" "C9"
$ws.Rows("10").RowHeight = 60

Set-Cell "A11" "Objectives:" "A9"
Set-Cell "B11" "TextB11" "B9"
Set-Cell "C11" "TextC11" "C9"
$ws.Rows("11").RowHeight = 60

Write-Host "done"
